# Remove the leftover "Card" rectangle/triangle/connector shapes that were
# cut from the bottom of the UI Component Class Diagram while the
# implementation was being reworked.
#
# Target shapes (matched by their PowerPoint shape Id, which is stable and
# unique within a slide):
#   id 92  - "Rectangle 11"          (text "Card")
#   id 105 - "Isosceles Triangle 102"
#   id 87  - "Connector: Elbow 86"
#   id 108 - "Connector: Elbow 107"
#   id 112 - "Connector: Elbow 111"

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$idsToDelete = @(92, 105, 87, 108, 112)

foreach ($targetId in $idsToDelete) {
    for ($i = $s.Shapes.Count; $i -ge 1; $i--) {
        $shp = $s.Shapes.Item($i)
        if ($shp.Id -eq $targetId) {
            $shp.Delete()
            break
        }
    }
}
